$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(' Oct 21 2020', ' Abu Dhabi', 'RCB won by 8 wickets (with 39 balls remaining)', 'Kolkata Knight Riders', 'Royal Challengers Bangalore', 'Shubman Gill ', '1', '6', '0', '0', '16.66'),
    @(' Oct 12 2020', ' Sharjah', 'RCB won by 82 runs', 'Kolkata Knight Riders', 'Royal Challengers Bangalore', 'Shubman Gill ', '34', '25', '3', '1', '136.00'),
    @(' Oct 16 2020', ' Abu Dhabi', 'Mumbai won by 8 wickets (with 19 balls remaining)', 'Kolkata Knight Riders', 'Mumbai Indians', 'Shubman Gill ', '21', '23', '2', '0', '91.30'),
    @(' Oct 3 2020', ' Sharjah', 'Capitals won by 18 runs', 'Kolkata Knight Riders', 'Delhi Capitals', 'Shubman Gill ', '28', '22', '2', '1', '127.27'),
    @(' Oct 24 2020', ' Abu Dhabi', 'KKR won by 59 runs', 'Kolkata Knight Riders', 'Delhi Capitals', 'Shubman Gill ', '9', '8', '2', '0', '112.50'),
    @(' Oct 18 2020', ' Abu Dhabi', 'Match tied (KKR won the one-over eliminator)', 'Kolkata Knight Riders', 'Sunrisers Hyderabad', 'Shubman Gill ', '36', '37', '5', '0', '97.29'),
    @(' Oct 7 2020', ' Abu Dhabi', 'KKR won by 10 runs', 'Kolkata Knight Riders', 'Chennai Super Kings', 'Shubman Gill ', '11', '12', '1', '0', '91.66'),
    @(' Oct 29 2020', ' Dubai (DSC)', 'Super Kings won by 6 wickets', 'Kolkata Knight Riders', 'Chennai Super Kings', 'Shubman Gill ', '26', '17', '4', '0', '152.94'),
    @(' Oct 26 2020', ' Sharjah', 'Kings XI won by 8 wickets (with 7 balls remaining)', 'Kolkata Knight Riders', 'Kings XI Punjab', 'Shubman Gill ', '57', '45', '3', '4', '126.66'),
    @(' Sep 26 2020', ' Abu Dhabi', 'KKR won by 7 wickets (with 12 balls remaining)', 'Kolkata Knight Riders', 'Sunrisers Hyderabad', 'Shubman Gill ', '70', '62', '5', '2', '112.90'),
    @(' Oct 10 2020', ' Abu Dhabi', 'KKR won by 2 runs', 'Kolkata Knight Riders', 'Kings XI Punjab', 'Shubman Gill ', '57', '47', '5', '0', '121.27'),
    @(' Nov 1 2020', ' Dubai (DSC)', 'KKR won by 60 runs', 'Kolkata Knight Riders', 'Rajasthan Royals', 'Shubman Gill ', '36', '24', '6', '0', '150.00'),
    @(' Sep 30 2020', ' Dubai (DSC)', 'KKR won by 37 runs', 'Kolkata Knight Riders', 'Rajasthan Royals', 'Shubman Gill ', '47', '34', '5', '1', '138.23'),
    @(' Sep 23 2020', ' Abu Dhabi', 'Mumbai won by 49 runs', 'Kolkata Knight Riders', 'Mumbai Indians', 'Shubman Gill ', '7', '11', '1', '0', '63.63'),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = "'" + $row[6]
    $ws.Cells.Item($r, 8).Value = "'" + $row[7]
    $ws.Cells.Item($r, 9).Value = "'" + $row[8]
    $ws.Cells.Item($r, 10).Value = "'" + $row[9]
    $ws.Cells.Item($r, 11).Value = "'" + $row[10]
}

$ws.Range("A1:K15").Select()
